$d = $word.ActiveDocument

$replacements = @(
    @{old="50×40=2000"; new="81×47=3807"},
    @{old="43×19=817"; new="46×28=1288"},
    @{old="73×15=1095"; new="47×17=799"},
    @{old="96×93=8928"; new="18×20=360"},
    @{old="91×90=8190"; new="56×81=4536"},
    @{old="98×67=6566"; new="45×67=3015"},
    @{old="36×36=1296"; new="73×34=2482"},
    @{old="87×48=4176"; new="52×17=884"},
    @{old="60×79=4740"; new="29×55=1595"},
    @{old="16×61=976"; new="85×27=2295"},
    @{old="68×13=884"; new="89×82=7298"},
    @{old="29×69=2001"; new="84×80=6720"},
    @{old="50×94=4700"; new="28×61=1708"},
    @{old="79×77=6083"; new="80×65=5200"},
    @{old="59×75=4425"; new="90×55=4950"},
    @{old="67×51=3417"; new="55×79=4345"},
    @{old="94×91=8554"; new="14×17=238"},
    @{old="93×15=1395"; new="32×42=1344"},
    @{old="70×79=5530"; new="11×86=946"},
    @{old="93×64=5952"; new="12×92=1104"},
    @{old="30×27=810"; new="87×43=3741"},
    @{old="69×66=4554"; new="23×28=644"},
    @{old="43×71=3053"; new="76×87=6612"},
    @{old="65×18=1170"; new="64×79=5056"},
    @{old="64×75=4800"; new="69×73=5037"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
